$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update model-year / MSRP figures for 2020 -> 2021 model year rollover ---
$updates = @(
    @{ Row = 2;  Year = 2021; Msrp = 42120 },
    @{ Row = 3;  Year = 2021; Msrp = 46590 },
    @{ Row = 4;  Year = 2021; Msrp = 44810 },
    @{ Row = 5;  Year = 2021; Msrp = 48765 },
    @{ Row = 6;  Year = 2021; Msrp = 45050 },
    @{ Row = 7;  Year = 2021; Msrp = 49520 },
    @{ Row = 8;  Year = 2021; Msrp = 47215 },
    @{ Row = 9;  Year = 2021; Msrp = 51130 },
    @{ Row = 53; Year = 2021; Msrp = 65875 },
    @{ Row = 54; Year = 2021; Msrp = 96675 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.Year
    $ws.Cells.Item($u.Row, 4).Value = $u.Msrp
}

# --- Append new "Black Line" trim rows (form deeplink test data) ---
$newRows = @(
    @{ Row = 95; Code = "9203SE"; Trim = "RC 300 F SPORT Black Line";     Year = 2021; Msrp = 48735 },
    @{ Row = 96; Code = "9207SE"; Trim = "RC 300 AWD F SPORT Black Line"; Year = 2021; Msrp = 50910 },
    @{ Row = 97; Code = "9213SE"; Trim = "RC 350 F SPORT Black Line";     Year = 2021; Msrp = 51665 },
    @{ Row = 98; Code = "9217SE"; Trim = "RC 350 AWD F SPORT Black Line"; Year = 2021; Msrp = 53275 }
)

# Reuse the number formats already used by the neighboring MSRP/DPHF column
# cells so the new rows pick up the same style instead of creating new ones.
$msrpNumberFormat = $ws.Cells.Item(94, 4).NumberFormat()
$dphfNumberFormat = $ws.Cells.Item(94, 5).NumberFormat()

# Populate column A (trim codes) for all new rows first, then column B (trim
# names), matching the order new strings are introduced in the workbook.
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Code
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Trim
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Year

    $ws.Cells.Item($r.Row, 4).NumberFormat = $msrpNumberFormat
    $ws.Cells.Item($r.Row, 4).Value = $r.Msrp

    $ws.Cells.Item($r.Row, 5).NumberFormat = $dphfNumberFormat
    $ws.Cells.Item($r.Row, 5).Value = 1025
}

# --- Update the view so the sheet opens scrolled/selected near the new rows ---
$win = $excel.ActiveWindow()
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("C55").Select()
